$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet 1: "Статистика по годам" - update numeric columns B, C, D, E for rows 2-9
$ws1.Cells.Item(2,2).Value = 44214
$ws1.Cells.Item(2,3).Value = 50240
$ws1.Cells.Item(2,4).Value = 35341
$ws1.Cells.Item(2,5).Value = 461
$ws1.Cells.Item(3,2).Value = 48536
$ws1.Cells.Item(3,3).Value = 48737
$ws1.Cells.Item(3,4).Value = 46657
$ws1.Cells.Item(3,5).Value = 613
$ws1.Cells.Item(4,2).Value = 44810
$ws1.Cells.Item(4,3).Value = 47379
$ws1.Cells.Item(4,4).Value = 31081
$ws1.Cells.Item(4,5).Value = 330
$ws1.Cells.Item(5,2).Value = 44648
$ws1.Cells.Item(5,3).Value = 48679
$ws1.Cells.Item(5,4).Value = 51686
$ws1.Cells.Item(5,5).Value = 630
$ws1.Cells.Item(6,2).Value = 46478
$ws1.Cells.Item(6,3).Value = 53581
$ws1.Cells.Item(6,4).Value = 77413
$ws1.Cells.Item(6,5).Value = 840
$ws1.Cells.Item(7,2).Value = 47924
$ws1.Cells.Item(7,3).Value = 61113
$ws1.Cells.Item(7,4).Value = 95147
$ws1.Cells.Item(7,5).Value = 945
$ws1.Cells.Item(8,2).Value = 53506
$ws1.Cells.Item(8,3).Value = 58409
$ws1.Cells.Item(8,4).Value = 129472
$ws1.Cells.Item(8,5).Value = 1011
$ws1.Cells.Item(9,2).Value = 49197
$ws1.Cells.Item(9,3).Value = 57885
$ws1.Cells.Item(9,4).Value = 141481
$ws1.Cells.Item(9,5).Value = 1250

$ws2 = $wb.Worksheets.Item(2)

# Sheet 2: "Статистика по городам" - update columns A,B (salary table) and D,E (share table)
$ws2.Cells.Item(2,1).Value = "Минск"
$ws2.Cells.Item(2,2).Value = 65722
$ws2.Cells.Item(2,4).Value = "Москва"
$ws2.Cells.Item(2,5).Formula = "'42.58%"
$ws2.Cells.Item(3,1).Value = "Москва"
$ws2.Cells.Item(3,2).Value = 59186
$ws2.Cells.Item(3,4).Value = "Санкт-Петербург"
$ws2.Cells.Item(3,5).Formula = "'12.51%"
$ws2.Cells.Item(4,1).Value = "Санкт-Петербург"
$ws2.Cells.Item(4,2).Value = 47008
$ws2.Cells.Item(4,4).Value = "Минск"
$ws2.Cells.Item(4,5).Formula = "'2.74%"
$ws2.Cells.Item(5,1).Value = "Новосибирск"
$ws2.Cells.Item(5,2).Value = 42055
$ws2.Cells.Item(5,4).Value = "Нижний Новгород"
$ws2.Cells.Item(5,5).Formula = "'2.57%"
$ws2.Cells.Item(6,1).Value = "Киев"
$ws2.Cells.Item(6,2).Value = 41172
$ws2.Cells.Item(6,4).Value = "Казань"
$ws2.Cells.Item(6,5).Formula = "'2.35%"
$ws2.Cells.Item(7,1).Value = "Екатеринбург"
$ws2.Cells.Item(7,2).Value = 40003
$ws2.Cells.Item(7,4).Value = "Новосибирск"
$ws2.Cells.Item(7,5).Formula = "'1.86%"
$ws2.Cells.Item(8,1).Value = "Алматы"
$ws2.Cells.Item(8,2).Value = 35147
$ws2.Cells.Item(8,4).Value = "Ростов-на-Дону"
$ws2.Cells.Item(8,5).Formula = "'1.81%"
$ws2.Cells.Item(9,1).Value = "Казань"
$ws2.Cells.Item(9,2).Value = 33339
$ws2.Cells.Item(9,4).Value = "Екатеринбург"
$ws2.Cells.Item(9,5).Formula = "'1.64%"
$ws2.Cells.Item(10,1).Value = "Нижний Новгород"
$ws2.Cells.Item(10,2).Value = 33023
$ws2.Cells.Item(10,4).Value = "Воронеж"
$ws2.Cells.Item(10,5).Formula = "'1.47%"
